# Weekly refresh of the Coliflor price series:
#  - Insert 2 new (blank) rows right before row 881, pushing the entire
#    existing 881..939 block down to 883..941 (brings the used range to
#    A1:R941).
#  - Fill the freshly-inserted 881/882 with the new weekly observation
#    (date 44931).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 881:939 down by two rows (A:R).
$ws.Range("A881:R882").Insert()

# New weekly data point, written into the now-freed rows 881/882.
# (Template columns A/B/C/E/F/G/H/N/O/Q/R are constant for every row in
# this block, so only the varying columns need to be (re)written.)
$ws.Range("A881").Value = 8
$ws.Range("B881").Value = "Terminal La Palmera de La Serena"
$ws.Range("C881").Value = "Coquimbo"
$ws.Range("D881").Value = 44931
$ws.Range("E881").Value = 4
$ws.Range("F881").Value = 100112008
$ws.Range("G881").Value = "Coliflor"
$ws.Range("H881").Value = "Sin especificar"
$ws.Range("I881").Value = "Primera"
$ws.Range("J881").Value = 2460
$ws.Range("K881").Value = 700
$ws.Range("L881").Value = 800
$ws.Range("M881").Value = 750
$ws.Range("N881").Value = "$/unidad"
$ws.Range("O881").Value = "Provincia del Elquí"
$ws.Range("P881").Value = 750
$ws.Range("Q881").Value = 1
$ws.Range("R881").Value = "Hortaliza"

$ws.Range("A882").Value = 8
$ws.Range("B882").Value = "Terminal La Palmera de La Serena"
$ws.Range("C882").Value = "Coquimbo"
$ws.Range("D882").Value = 44931
$ws.Range("E882").Value = 4
$ws.Range("F882").Value = 100112008
$ws.Range("G882").Value = "Coliflor"
$ws.Range("H882").Value = "Sin especificar"
$ws.Range("I882").Value = "Segunda"
$ws.Range("J882").Value = 1480
$ws.Range("K882").Value = 500
$ws.Range("L882").Value = 600
$ws.Range("M882").Value = 550
$ws.Range("N882").Value = "$/unidad"
$ws.Range("O882").Value = "Provincia del Elquí"
$ws.Range("P882").Value = 550
$ws.Range("Q882").Value = 1
$ws.Range("R882").Value = "Hortaliza"
